# Updates the "Report" sheet's two correlation tables:
#  - the numeric correlation matrix for IdTrx/Amount/Datekey (rows 46-48)
#  - the "top correlated categorical values" list (rows 57-85), whose
#    labels and correlation scores were refreshed after re-running the
#    report against newer data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric correlation matrix (rows 46-48) ---
$ws.Range("C46").Value = 0.04375210694513654
$ws.Range("D46").Value = 0.9823871556748695
$ws.Range("B47").Value = 0.04375210694513654
$ws.Range("D47").Value = 0.05271375803118115
$ws.Range("B48").Value = 0.9823871556748695
$ws.Range("C48").Value = 0.05271375803118115

# --- Top correlated categorical values (rows 57-85): label + score ---
$ws.Range("A57").Value = "Instruction Date_2020-02-27"
$ws.Range("B57").Value = 0.92
$ws.Range("A58").Value = "BIC FI Initiation_CHASUS33FXR"
$ws.Range("B58").Value = 0.21
$ws.Range("A59").Value = "Instruction Date_2014-11-20"
$ws.Range("B59").Value = 0.21
$ws.Range("A60").Value = "BIC FI Initiation_ROYCGB2L"
$ws.Range("B60").Value = 0.16
$ws.Range("A61").Value = "Instruction Date_2015-10-21"
$ws.Range("B61").Value = 0.16
$ws.Range("A62").Value = "BIC Sender_ROYCCAT2"
$ws.Range("B62").Value = 0.14
$ws.Range("A63").Value = "Instruction Date_2018-07-15"
$ws.Range("B63").Value = 0.12
$ws.Range("A64").Value = "BIC FI Initiation_ROYCCAT2"
$ws.Range("B64").Value = 0.11
$ws.Range("A65").Value = "Instruction Date_2016-09-17"
$ws.Range("B65").Value = 0.11
$ws.Range("A66").Value = "BIC FI Destination_BNDCCAMMINT"
$ws.Range("B66").Value = 0.11
$ws.Range("A67").Value = "ind_lvts_Oui"
$ws.Range("B67").Value = 0.11
$ws.Range("A68").Value = "Instruction Date_2020-02-26"
$ws.Range("B68").Value = 0.1
$ws.Range("A69").Value = "Country Sender_CA"
$ws.Range("B69").Value = 0.09
$ws.Range("A70").Value = "Country FI Initiation_CA"
$ws.Range("B70").Value = 0.08
$ws.Range("A71").Value = "BIC FI Initiation_BNDCCAMM"
$ws.Range("B71").Value = 0.08
$ws.Range("A72").Value = "is_sender FI Initiation_Oui"
$ws.Range("B72").Value = 0.08
$ws.Range("A73").Value = "Channel_UNKNOWN"
$ws.Range("B73").Value = 0.07
$ws.Range("A74").Value = "BIC Sender_BCANCAW2"
$ws.Range("B74").Value = 0.07
$ws.Range("A75").Value = "BIC FI Initiation_BCANCAW2"
$ws.Range("B75").Value = 0.07
$ws.Range("A76").Value = "Instruction Date_2014-12-27"
$ws.Range("B76").Value = 0.07
$ws.Range("A77").Value = "Instruction Date_2018-01-26"
$ws.Range("B77").Value = 0.07
$ws.Range("A78").Value = "Solution_Product_Direct_Debit"
$ws.Range("B78").Value = 0.07
$ws.Range("A79").Value = "Instruction_Withdrawal"
$ws.Range("B79").Value = 0.07
$ws.Range("A80").Value = "Datekey"
$ws.Range("A81").Value = "Instruction Date_2017-07-11"
$ws.Range("A82").Value = "Country FI Initiation_GB"
$ws.Range("A83").Value = "Country FI Initiation_US"
$ws.Range("B83").Value = 0.04
$ws.Range("A84").Value = "IdTrx"
$ws.Range("B84").Value = 0.04
$ws.Range("A85").Value = "Instruction Date_2020-01-26"
